# Seed database with all priority levels
#
# Rename the "priority_level" sheet to "priority" and update its header
# row from priority_level_id/priority_level to priority_id/name, then
# leave that sheet active/selected (matching the author's final view
# state before saving).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("priority_level")
$ws.Name = "priority"

$ws.Range("A1").Value = "priority_id"
$ws.Range("B1").Value = "name"

$ws.Activate()
$ws.Range("C14").Select()
